# Cycle-3 data refresh: raw annotation interval duration (E) and the
# matching Aggregates interval length (C) were re-measured at a uniform
# 300s window, replacing the assorted per-flow durations. Downstream
# throughput/flow-rate formulas (J/K on Raw_Annotations, F/H on
# Aggregates) recalc automatically.

$wb = $excel.ActiveWorkbook

$wsRaw = $wb.Worksheets.Item("Raw_Annotations")
$wsAgg = $wb.Worksheets.Item("Aggregates")

# --- Raw_Annotations: column E (interval seconds) -> 300 for all data rows
$wsRaw.Range("E2:E31").Value = 300

# --- Aggregates: column C (interval seconds) -> 300 for all data rows
$wsAgg.Range("C2:C6").Value = 300

# --- Selection bookkeeping matching the authored session: Aggregates was
# selected first (C2:C6), then the user switched to Raw_Annotations and
# selected E2:E31, leaving Raw_Annotations as the final active sheet/tab.
$wsAgg.Activate() | Out-Null
$wsAgg.Range("C2:C6").Select() | Out-Null

$wsRaw.Activate() | Out-Null
$wsRaw.Range("E2:E31").Select() | Out-Null

Write-Output "Updated interval durations to 300s across Raw_Annotations (E2:E31) and Aggregates (C2:C6)."
